$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet Q1: add "Claim/Challenge vs Status-quo" callout next to the
# existing Ho/Ha block (I4:I5), with a thin right border on I4 only
# (echoing the existing boxed-label look used elsewhere in the sheet).
# ---------------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item("Q1")
$wsQ1.Range("I4").Value = "Statu quo "
$wsQ1.Range("I4").Borders.Item(10).LineStyle = 1
$wsQ1.Range("I5").Value = "Claim or challenge"

# ---------------------------------------------------------------------------
# Sheet Q2: tighten the Ho/Ha proportion thresholds from the boundary value
# 0.58 to 0.57999 (so Ho/Ha partition strictly), and add the companion
# Claim/Status-quo + inequality table (I4:J5, J7).
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("Q2")
$wsQ2.Range("H4").Value = "Pi <= 0.57999"
$wsQ2.Range("H5").Value = "Pi > 0.57999"

$wsQ2.Range("I4").Value = "Status quo "
$wsQ2.Range("I4").Borders.Item(7).LineStyle = 1
$wsQ2.Range("I4").Borders.Item(10).LineStyle = 1

$wsQ2.Range("I5").Value = "Claim "
$wsQ2.Range("I5").Borders.Item(7).LineStyle = 1
$wsQ2.Range("I5").Borders.Item(10).LineStyle = 1

# J4/J5/J7 reuse the same "boxed label" look already used by cells like H7
# (font/fill/border) -- copy that formatting across rather than rebuilding it.
$wsQ2.Range("H7").Copy() | Out-Null
$wsQ2.Range("J4").PasteSpecial(-4122) | Out-Null
$wsQ2.Range("J5").PasteSpecial(-4122) | Out-Null
$wsQ2.Range("J7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsQ2.Range("J4").Value = "Pi >= 0.58"
$wsQ2.Range("J5").Value = "Pi < 0.58"
$wsQ2.Range("J7").Value = "LeftTail test "

# ---------------------------------------------------------------------------
# Sheet Q3: note the left-tail hypotheses next to the existing statistics.
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("Q3")
$wsQ3.Range("I3").Value = "mud <=5 5"
$wsQ3.Range("I4").Value = "mud>=5"

# ---------------------------------------------------------------------------
# Sheet Q4: add pooled-standard-deviation notes (L7:L8), fix the upper
# critical-t-value formulas to reference the degrees of freedom cell B9
# (not the T statistic in B10), and flip the rejection-rule comparisons to
# use a right-tail ">" test instead of a two-tail "<" test.
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("Q4")
$wsQ4.Range("L7").Value = "sp = n1-2 "
$wsQ4.Range("L8").Value = "sp * sqrt(1/n1 + 1/n2)"

$wsQ4.Range("J10").Formula = '=T.INV(H10,$B$9)'
$wsQ4.Range("J11").Formula = '=T.INV(H11,$B$9)'
$wsQ4.Range("J12").Formula = '=T.INV(H12,$B$9)'

$wsQ4.Range("C15").Formula = '=IF(ABS(I10)>$C$7,"Reject Null Hypothesis","Fail to reject Null Hypothesis")'
$wsQ4.Range("C16").Formula = '=IF(ABS(I11)>$C$7,"Reject Null Hypothesis","Fail to reject Null Hypothesis")'
$wsQ4.Range("C17").Formula = '=IF(ABS(I12)>$C$7,"Reject Null Hypothesis","Fail to reject Null Hypothesis")'

$wb.Save()
